$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1569.5714
$ws.Range("I38").Value = 30.333334
$ws.Range("K38").Value = 91.00000199999999
$ws.Range("M38").Value = 280.999998
$ws.Range("H86").Value = 83336170
$ws.Range("I86").Value = 100003090
$ws.Range("J86").Value = 1570
$ws.Range("K86").Value = 100003090
$ws.Range("L86").Value = 1570
$ws.Range("M86").Value = -100001967
$ws.Range("N86").Value = -3816
$ws.Range("H89").Value = 83336170
$ws.Range("I89").Value = 100003090
$ws.Range("J89").Value = 1570
$ws.Range("K89").Value = 500015450
$ws.Range("L89").Value = 7850
$ws.Range("M89").Value = -500009834
$ws.Range("N89").Value = -19082
$ws.Range("H132").Value = 2855
$ws.Range("I132").Value = 2841.238
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 8523.714
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -5993.714
$ws.Range("N132").Value = -14058.5
$ws.Range("H137").Value = 2569943.8
$ws.Range("I137").Value = 4676.0835
$ws.Range("J137").Value = 6674372
$ws.Range("K137").Value = 14028.2505
$ws.Range("L137").Value = 20023116
$ws.Range("M137").Value = -11478.2505
$ws.Range("N137").Value = -20028216
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 9994.666999999999
$ws.Range("J36").Value = 9993.333000000001
$ws.Range("L36").Value = 9993.333000000001
$ws.Range("N36").Value = -10685.333
$ws.Range("H37").Value = 53998.75
$ws.Range("J37").Value = 56999
$ws.Range("L37").Value = 56999
$ws.Range("N37").Value = -57545
$ws.Range("H44").Value = 66682
$ws.Range("J44").Value = 99999
$ws.Range("L44").Value = 99999
$ws.Range("N44").Value = -100975
$ws.Range("H45").Value = 49723.668
$ws.Range("I45").Value = 73014.14
$ws.Range("K45").Value = 73014.14
$ws.Range("M45").Value = -72637.14
$ws.Range("H55").Value = 52499.5
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 99999
$ws.Range("N55").Value = -100629
$ws.Range("H63").Value = 23922.158
$ws.Range("I63").Value = 3164
$ws.Range("K63").Value = 3164
$ws.Range("M63").Value = -2478
$ws.Range("H66").Value = 23922.158
$ws.Range("I66").Value = 3164
$ws.Range("K66").Value = 15820
$ws.Range("M66").Value = -12388
$ws.Range("H74").Value = 619849.6
$ws.Range("I74").Value = 1037.4736
$ws.Range("J74").Value = 1311463.2
$ws.Range("K74").Value = 1037.4736
$ws.Range("L74").Value = 1311463.2
$ws.Range("M74").Value = -163.4736
$ws.Range("N74").Value = -1313211.2
$ws.Range("H77").Value = 619849.6
$ws.Range("I77").Value = 1037.4736
$ws.Range("J77").Value = 1311463.2
$ws.Range("K77").Value = 5187.368
$ws.Range("L77").Value = 6557316
$ws.Range("M77").Value = -819.3680000000004
$ws.Range("N77").Value = -6566052
$ws.Range("H80").Value = 51999.5
$ws.Range("J80").Value = 51999.5
$ws.Range("L80").Value = 51999.5
$ws.Range("N80").Value = -53995.5
$ws.Range("H83").Value = 51999.5
$ws.Range("J83").Value = 51999.5
$ws.Range("L83").Value = 155998.5
$ws.Range("N83").Value = -165982.5
$ws.Range("H88").Value = 1027
$ws.Range("J88").Value = 346.5
$ws.Range("L88").Value = 346.5
$ws.Range("N88").Value = -1158.5
$ws.Range("H91").Value = 1027
$ws.Range("J91").Value = 346.5
$ws.Range("L91").Value = 346.5
$ws.Range("N91").Value = -3154.5
$ws.Range("H102").Value = 2066.25
$ws.Range("I102").Value = 2110.5652
$ws.Range("J102").Value = 1047
$ws.Range("K102").Value = 2110.5652
$ws.Range("L102").Value = 1047
$ws.Range("M102").Value = -488.5652
$ws.Range("N102").Value = -4291
$ws.Range("H132").Value = 3119.9412
$ws.Range("I132").Value = 3367.7
$ws.Range("K132").Value = 10103.1
$ws.Range("M132").Value = -7573.099999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 29034340
$ws.Range("I134").Value = 1701.56
$ws.Range("K134").Value = 5104.68
$ws.Range("M134").Value = -2569.68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2721.4883
$ws.Range("I31").Value = 3417.2942
$ws.Range("K31").Value = 3417.2942
$ws.Range("M31").Value = -3122.2942
$ws.Range("H34").Value = 2721.4883
$ws.Range("I34").Value = 3417.2942
$ws.Range("K34").Value = 3417.2942
$ws.Range("M34").Value = -3215.2942
$ws.Range("H62").Value = 3594.2354
$ws.Range("I62").Value = 3918.5557
$ws.Range("J62").Value = 3229.375
$ws.Range("K62").Value = 3918.5557
$ws.Range("L62").Value = 3229.375
$ws.Range("M62").Value = -3294.5557
$ws.Range("N62").Value = -4477.375
$ws.Range("H65").Value = 3594.2354
$ws.Range("I65").Value = 3918.5557
$ws.Range("J65").Value = 3229.375
$ws.Range("K65").Value = 19592.7785
$ws.Range("L65").Value = 16146.875
$ws.Range("M65").Value = -16472.7785
$ws.Range("N65").Value = -22386.875
$ws.Range("H107").Value = 1596.6471
$ws.Range("I107").Value = 1680.9286
$ws.Range("K107").Value = 1680.9286
$ws.Range("M107").Value = 239.0714
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1400.3334
$ws.Range("I17").Value = 500.5
$ws.Range("K17").Value = 1501.5
$ws.Range("M17").Value = -1332.5
$ws.Range("H39").Value = 5248.5
$ws.Range("J39").Value = 6266.875
$ws.Range("L39").Value = 18800.625
$ws.Range("N39").Value = -19388.625
$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 8000
$ws.Range("L50").Value = 24000
$ws.Range("N50").Value = -24962
$ws.Range("H53").Value = 8000
$ws.Range("J53").Value = 8000
$ws.Range("L53").Value = 24000
$ws.Range("N53").Value = -24962
$ws.Range("H107").Value = 1326.2307
$ws.Range("J107").Value = 1811.125
$ws.Range("L107").Value = 5433.375
$ws.Range("N107").Value = -9273.375
$ws.Range("H114").Value = 1395.5333
$ws.Range("I114").Value = 1099
$ws.Range("J114").Value = 1593.2222
$ws.Range("K114").Value = 3297
$ws.Range("L114").Value = 4779.6666
$ws.Range("M114").Value = -43
$ws.Range("N114").Value = -11287.6666
$ws.Range("H133").Value = 5250
$ws.Range("I133").Value = 5250
$ws.Range("K133").Value = 15750
$ws.Range("M133").Value = -10690
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 80000
$ws.Range("I108").Value = 80000
$ws.Range("K108").Value = 80000
$ws.Range("M108").Value = -76160
$ws.Range("H132").Value = 2851213.5
$ws.Range("I132").Value = 1574.25
$ws.Range("J132").Value = 5130925
$ws.Range("K132").Value = 4722.75
$ws.Range("L132").Value = 15392775
$ws.Range("M132").Value = -2192.75
$ws.Range("N132").Value = -15397835
$ws.Range("H136").Value = 27999
$ws.Range("J136").Value = 27999
$ws.Range("L136").Value = 83997
$ws.Range("N136").Value = -89097
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4034.125
$ws.Range("I22").Value = 1138.2222
$ws.Range("J22").Value = 5771.6665
$ws.Range("K22").Value = 1138.2222
$ws.Range("L22").Value = 5771.6665
$ws.Range("M22").Value = -843.2221999999999
$ws.Range("N22").Value = -6361.6665
$ws.Range("H27").Value = 4034.125
$ws.Range("I27").Value = 1138.2222
$ws.Range("J27").Value = 5771.6665
$ws.Range("K27").Value = 1138.2222
$ws.Range("L27").Value = 5771.6665
$ws.Range("M27").Value = -1031.2222
$ws.Range("N27").Value = -5985.6665
$ws.Range("H59").Value = 85500
$ws.Range("J59").Value = 85500
$ws.Range("L59").Value = 85500
$ws.Range("N59").Value = -86808
$ws.Range("H82").Value = 2620.25
$ws.Range("I82").Value = 1494
$ws.Range("J82").Value = 5999
$ws.Range("K82").Value = 1494
$ws.Range("L82").Value = 5999
$ws.Range("M82").Value = -1133
$ws.Range("N82").Value = -6721
$ws.Range("H85").Value = 2620.25
$ws.Range("I85").Value = 1494
$ws.Range("J85").Value = 5999
$ws.Range("K85").Value = 1494
$ws.Range("L85").Value = 5999
$ws.Range("M85").Value = -246
$ws.Range("N85").Value = -8495
$ws.Range("H132").Value = 3413.7222
$ws.Range("I132").Value = 3335.7
$ws.Range("J132").Value = 3511.25
$ws.Range("K132").Value = 10007.1
$ws.Range("L132").Value = 10533.75
$ws.Range("M132").Value = -7477.099999999999
$ws.Range("N132").Value = -15593.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27999.8
$ws.Range("H107").Value = 4762369
$ws.Range("I107").Value = 263.33334
$ws.Range("J107").Value = 9524474
$ws.Range("K107").Value = 790.0000200000001
$ws.Range("L107").Value = 28573422
$ws.Range("M107").Value = 1129.99998
$ws.Range("N107").Value = -28577262
$ws.Range("H132").Value = 2041.1351
$ws.Range("I132").Value = 1522.6666
$ws.Range("J132").Value = 3441
$ws.Range("K132").Value = 4567.9998
$ws.Range("L132").Value = 10323
$ws.Range("M132").Value = -2037.9998
$ws.Range("N132").Value = -15383
